$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The statistics table is being re-ordered (Медицина/Математика/Лингвистика/Физика
# -> Математика/Лингвистика/Медицина/Физика) and the "average score" column (B) is
# now written as text (BigDecimal-backed) instead of a floating point number for the
# two rows whose average is a whole number (0 and 5).

# Row 2: Математика
$ws.Range("A2").Value = "Математика"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = 0.0
$ws.Range("D2").Value = 1.0
$ws.Range("E2").Value = "Казанский Университет Вычислений;"

# Row 3: Лингвистика
$ws.Range("A3").Value = "Лингвистика"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = 0.0
$ws.Range("D3").Value = 1.0
$ws.Range("E3").Value = "Воронежский Литературно-Переводческий Университет;"

# Row 4: Медицина
$ws.Range("A4").Value = "Медицина"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "5"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = 3.0
$ws.Range("D4").Value = 3.0
$ws.Range("E4").Value = "Московский Государственный Медицинский Университет;Тамбовский Университет Медицины;Самарский Медицинский Институт;"

# Row 5: Физика (stays last, but its average score is now text too)
$ws.Range("A5").Value = "Физика"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "5"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = 8.0
$ws.Range("D5").Value = 2.0
$ws.Range("E5").Value = "Московский Выдуманный Университет;Московский Придуманный Институт;"
